# Update gh-pages to output generated at 456a3b4
# Applies to workbook "苏州-漫展信息.xlsx":
#  - Sheet "展览" (exhibitions) and "全部类型" (all types): remove three
#    cancelled/irrelevant events, renumber the index column, and bump a
#    handful of "想去人数" (interested count) figures.
#  - Sheet "演出" (performances): bump two "想去人数" figures.

$wb = $excel.ActiveWorkbook

$namesToDelete = @(
    "常熟·cc动漫游戏嘉年华",
    "张家港· 喵喵漫国潮动漫节",
    "苏州·漫遇引力动漫游戏展"
)

function Remove-EventRows($ws, $names) {
    foreach ($name in $names) {
        $found = $ws.Cells.Find($name)
        if ($found -ne $null) {
            $found.EntireRow.Delete()
        }
    }
}

function Renumber-IndexColumn($ws) {
    $lastRow = $ws.Cells(1, 1).End(-4121).Row
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }
}

function Set-WantCount($ws, $updates) {
    foreach ($name in $updates.Keys) {
        $found = $ws.Cells.Find($name)
        if ($found -ne $null) {
            $ws.Cells.Item($found.Row, 6).Value = $updates[$name]
        }
    }
}

# ---------------------------------------------------------------------
# Sheet "展览"
# ---------------------------------------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")
Remove-EventRows $wsExhibit $namesToDelete
Renumber-IndexColumn $wsExhibit

$exhibitUpdates = @{
    "苏州·恋与深空only 同人周边套餐 "              = 761
    "苏州·无限次元夜场"                            = 109
    "苏州·国内知名配音演员吕书君@阿君归来 签售会"   = 9
    "苏州·漫语堂动漫嘉年华"                        = 107
    "苏州·女神异闻录only同人展"                    = 417
    "张家港·META萌圆饿了"                          = 129
    "苏州·COME IN JOY 动漫品牌国潮文化节"          = 11249
    "苏州·星部落&青铜树动漫嘉年华"                 = 5345
}
Set-WantCount $wsExhibit $exhibitUpdates

# ---------------------------------------------------------------------
# Sheet "演出"
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$showUpdates = @{
    "苏州·乐队番同人only live Band Set二次元乐队拼盘" = 93
    "苏州·维也纳皇家交响乐团2025新年音乐会"           = 5
}
Set-WantCount $wsShow $showUpdates

# ---------------------------------------------------------------------
# Sheet "全部类型"
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
Remove-EventRows $wsAll $namesToDelete
Renumber-IndexColumn $wsAll

$allUpdates = @{
    "苏州·恋与深空only 同人周边套餐 "              = 761
    "苏州·明日方舟同人only （聚会）"               = 52
    "苏州·乐队番同人only live Band Set二次元乐队拼盘" = 93
    "苏州·无限次元夜场"                            = 109
    "苏州·国内知名配音演员吕书君@阿君归来 签售会"   = 9
    "苏州·漫语堂动漫嘉年华"                        = 107
    "苏州·女神异闻录only同人展"                    = 417
    "张家港·META萌圆饿了"                          = 129
    "苏州·COME IN JOY 动漫品牌国潮文化节"          = 11249
    "苏州·维也纳皇家交响乐团2025新年音乐会"           = 5
    "苏州·星部落&青铜树动漫嘉年华"                 = 5345
}
Set-WantCount $wsAll $allUpdates
